$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.993.79'
$ws.Range('E2').Value = '  -0.55%  '
$ws.Range('D3').Value = '1.828.12'
$ws.Range('E3').Value = '  +0.00%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '311.56'
$ws.Range('E5').Value = '  -0.56%  '
$ws.Range('E6').Value = '  -0.11%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4653'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3700'
$ws.Range('E8').Value = '  +1.35%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07355'
$ws.Range('E9').Value = '  -0.57%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8729'
$ws.Range('E10').Value = '  -0.89%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07880'
$ws.Range('E11').Value = '  +7.53%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '19.91'
$ws.Range('E12').Value = '  -2.15%  '
$ws.Range('D13').Value = '1.771.90'
$ws.Range('E13').Value = '  -8.17%  '
$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.596'
$ws.Range('E14').Value = '  +1.16%  '
$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.355'
$ws.Range('E15').Value = '  -0.49%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '91.93'
$ws.Range('E16').Value = '  -1.49%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.010'
$ws.Range('E17').Value = '  +0.25%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008869'
$ws.Range('E18').Value = '  +1.93%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.008'
$ws.Range('E19').Value = '  -0.11%  '
$ws.Range('E20').Value = '  +0.36%  '
$ws.Range('D21').Value = '27.190.09'
$ws.Range('E21').Value = '  -1.72%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.154'
$ws.Range('E22').Value = '  -1.66%  '
$ws.Range('E23').Value = '  -0.29%  '
$ws.Range('D24').Value = '2.020.71'
$ws.Range('E24').Value = '  -3.71%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '152.59'
$ws.Range('E25').Value = '  +0.57%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.832'
$ws.Range('E26').Value = '  -2.63%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.19'
$ws.Range('E27').Value = '  -1.80%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.095'
$ws.Range('E28').Value = '  -2.19%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.125'
$ws.Range('E29').Value = '  -1.14%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '115.48'
$ws.Range('E30').Value = '  -0.63%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08871'
$ws.Range('E31').Value = '  -0.76%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.985'
$ws.Range('E32').Value = '  +1.63%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.7271'
$ws.Range('E33').Value = '  -2.07%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.440'
$ws.Range('E34').Value = '  -1.56%  '
$ws.Range('E35').Value = '  -2.78%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.497'
$ws.Range('E36').Value = '  +3.42%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.078'
$ws.Range('E37').Value = '  -1.14%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01949'
$ws.Range('E38').Value = '  +0.01%  '
$ws.Range('E39').Value = '  -1.20%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '7.285'
$ws.Range('E40').Value = '  +1.16%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.929'
$ws.Range('E41').Value = '  -0.15%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5188'
$ws.Range('E42').Value = '  -1.23%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8653'
$ws.Range('E43').Value = '  -14.33%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.1624'
$ws.Range('E44').Value = '  -1.15%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.209'
$ws.Range('E45').Value = '  -2.07%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4836'
$ws.Range('E46').Value = '  -0.65%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.009'
$ws.Range('E47').Value = '  -0.13%  '
$ws.Range('E48').Value = '  -2.08%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '102.76'
$ws.Range('E49').Value = '  -1.44%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.623'
$ws.Range('E50').Value = '  -1.69%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06227'
$ws.Range('E51').Value = '  -1.13%  '
